$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 into the new
# header cells I1 and J1 before setting their values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row (row 1) - new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-8 for column I (I0) and column J (IF)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 6

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 6

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 6

$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 5

$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 6
